# Update ticker-style market names in columns A/B to Korean coin names
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "리플"
$ws.Range("B2").Value = "스텔라루멘"

$ws.Range("A3").Value = "쿼크체인"
$ws.Range("B3").Value = "아이오에스티"

# Add a new row of data (price-band coupling example)
$ws.Range("A4").Value = "리스크"
$ws.Range("B4").Value = "오미세고"
$ws.Range("C4").Value = "rare"
$ws.Range("D4").Value = "1000~3000원"
$ws.Range("E4").Value = "가격대 커플링"

# Widen columns A and B now that they hold longer Korean text
$ws.Columns.Item(1).ColumnWidth = 15.142857142857142
$ws.Columns.Item(2).ColumnWidth = 20.857142857142854

# Move the active selection to the newly added cell
$ws.Range("E4").Select()
